# "Conversion de Matriz a arreglo lineal" - adds a second timing table
# (columns F & O, rows 38-80) below the existing one, plus two summary
# formulas in row 67 (F) and row 81 (O/P), mirroring the existing
# "Tot"/"Duracion" pattern used earlier in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- new data block: column F ("F..") and column O ("nseg"-like series) ---
$ws.Range("O38").Value = 3
$ws.Range("F39").Value = 3
$ws.Range("O39").Value = 2
$ws.Range("F40").Value = 2
$ws.Range("O40").Value = 1
$ws.Range("F41").Value = 1
$ws.Range("O41").Value = 3
$ws.Range("F42").Value = 1
$ws.Range("O42").Value = 1
$ws.Range("F43").Value = 3
$ws.Range("O43").Value = 3
$ws.Range("F44").Value = 2
$ws.Range("O44").Value = 1
$ws.Range("F45").Value = 1
$ws.Range("O45").Value = 2
$ws.Range("F46").Value = 2
$ws.Range("O46").Value = 3
$ws.Range("F47").Value = 1
$ws.Range("O47").Value = 2
$ws.Range("F48").Value = 1
$ws.Range("O48").Value = 3
$ws.Range("F49").Value = 3
$ws.Range("O49").Value = 2
$ws.Range("F50").Value = 1
$ws.Range("O50").Value = 2
$ws.Range("F51").Value = 2
$ws.Range("O51").Value = 1
$ws.Range("F52").Value = 2
$ws.Range("O52").Value = 3
$ws.Range("F53").Value = 3
$ws.Range("O53").Value = 1
$ws.Range("F54").Value = 1
$ws.Range("O54").Value = 2
$ws.Range("F55").Value = 3
$ws.Range("O55").Value = 1
$ws.Range("F56").Value = 1
$ws.Range("O56").Value = 1
$ws.Range("F57").Value = 4
$ws.Range("O57").Value = 3
$ws.Range("F58").Value = 3
$ws.Range("O58").Value = 1
$ws.Range("F59").Value = 4
$ws.Range("O59").Value = 2
$ws.Range("F60").Value = 3
$ws.Range("O60").Value = 4
$ws.Range("F61").Value = 5
$ws.Range("O61").Value = 3
$ws.Range("F62").Value = 3
$ws.Range("O62").Value = 3
$ws.Range("F63").Value = 2
$ws.Range("O63").Value = 4
$ws.Range("F64").Value = 2
$ws.Range("O64").Value = 3
$ws.Range("F65").Value = 3
$ws.Range("O65").Value = 0
$ws.Range("O66").Value = 0
$ws.Range("O67").Value = 0
$ws.Range("O68").Value = 4
$ws.Range("O69").Value = 2
$ws.Range("O70").Value = 3
$ws.Range("O71").Value = 3
$ws.Range("O72").Value = 3
$ws.Range("O73").Value = 3
$ws.Range("O74").Value = 0
$ws.Range("O75").Value = 4
$ws.Range("O76").Value = 2
$ws.Range("O77").Value = 3
$ws.Range("O78").Value = 2
$ws.Range("O79").Value = 3
$ws.Range("O80").Value = 3

# --- summary formulas ---
$ws.Range("F67").Formula = "=25*32*8*SUM(F39:F65)"
$ws.Range("O81").Formula = "=SUM(O38:O80)"
$ws.Range("P81").Formula = "=32*8*O81*25"

# --- recalc + view state (active cell / scroll position) ---
$excel.Calculate()
$ws.Range("P81").Select()
$excel.ActiveWindow.ScrollRow = 61
$excel.ActiveWindow.ScrollColumn = 3
